$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "Tomos" (volumes owned) counts for a handful of series.
# The dependent "Valor total" (G column, D*E) and the summary rows
# (96/98) recalc automatically as formulas.
$ws.Range("D2").Value = 20
$ws.Range("D7").Value = 11
$ws.Range("D11").Value = 11
$ws.Range("D13").Value = 6
$ws.Range("D25").Value = 2

# Move the current selection / scroll position, matching the
# author's last on-screen view before saving.
$ws.Range("A22").Select()
$excel.ActiveWindow.ScrollRow = 22
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("D26").Select()
